$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 (pushes the existing rows 23-33 down to 24-34,
# extending the data block with a new weekly price record).
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly record.
$ws.Range("A23").Value = 7
$ws.Range("B23").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C23").Value = "Ñuble"
$ws.Range("D23").Value = 44830
$ws.Range("E23").Value = 16
$ws.Range("F23").Value = 100112043
$ws.Range("G23").Value = "Pepino dulce"
$ws.Range("H23").Value = "Cultivar IV Región"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 60
$ws.Range("K23").Value = 17000
$ws.Range("L23").Value = 17000
$ws.Range("M23").Value = 17000
$ws.Range("N23").Value = "$/bandeja 18 kilos"
$ws.Range("O23").Value = "Provincia de Limarí"
$ws.Range("P23").Value = 944
$ws.Range("Q23").Value = 18
$ws.Range("R23").Value = "Hortaliza"
